$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(40).Insert() | Out-Null
$ws.Range("A40").Value = "topography_oasis"
$ws.Range("B40").Value = "Oasis"
$ws.Range("A40").Select() | Out-Null
